$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "23.977.77"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.653.00"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +1.95%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "308.78"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.30%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9990"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3903"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("E8").Value = "  -0.37%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "51.28"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.39%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.356"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9987"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.06%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08446"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "23.95"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +1.01%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.862"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +4.04%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001316"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +2.73%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.653.05"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "94.46"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.72%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06975"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.61%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "19.79"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.890"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("E22").Value = "  +0.07%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.62"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "23.974.35"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.476"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.70%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "3.031"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.06%  "
$ws.Range("E27").Value = "  -0.96%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "152.66"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.48%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.443"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.26%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "139.22"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.47%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.754"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -1.38%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.485"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.29%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.835.40"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.90%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.035"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +5.08%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.08057"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("E36").Value = "  +2.73%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.709"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  +4.74%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.2683"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.09124"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.37%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.7561"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.64%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "13.47"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.78%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.435"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.24"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6941"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.50%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.451"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -1.02%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.071"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.12%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.9984"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.08317"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.80%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "134.34"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.46%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "1.223"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.20%  "
